$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Periodo Mora" period for every worker row (2507 -> 2508).
#    All rows 16-31 share the same value, so update the whole column at once.
$ws.Range("E16:E31").Value = "2508"

# 2. Update aggregate "VALOR MORA" total.
$ws.Range("E11").Value = 868220

# 3. Update "Cant. Trabajadores" count (one worker is being removed below).
$ws.Range("C13").Value = 15

# 4. Update "Salario Basico" for JAVIER ENRIQUE FORTICH PRENS (row 17).
$ws.Range("G17").Value = 1500000

# 5. Remove the last worker row (EMIRO DE JESUS LARA TORRES / 73086885).
#    First copy its formatting (thicker bottom border) up onto the row that
#    will become the new last row, then delete the row so everything below
#    shifts up.
$ws.Range("B31:J31").Copy()
$ws.Range("B30:J30").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows("31:31").Delete()
